$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute("September 19, 2025", $false, $false, $false, $false, $false, `
                                 $true, 1, $false, "September 21, 2025", 2)

# ------------------------------------------------------------------
# 2. Split the mailing-address line into two paragraphs:
#    "969 Story Road, San Jose CA 95122"
#      -> "969 Story Road"
#         "San Jose, CA 95122"
#    Only the first occurrence (the addressee block near the top of
#    the letter) is affected; the "PROPERTY ADDRESS:" block further
#    down is left untouched.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "969 Story Road, San Jose CA 95122") {
        $r = $p.Range
        $r.Text = "969 Story Road"
        $r.InsertParagraphAfter() | Out-Null
        $d.Paragraphs($i + 1).Range.Text = "San Jose, CA 95122"
        break
    }
}

# ------------------------------------------------------------------
# 3. Remove the empty "No Spacing" paragraph that immediately follows
#    the "... Board of Directors" signature line.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -match "Board of Directors$") {
        $next = $d.Paragraphs($i + 1)
        $nextText = $next.Range.Text.TrimEnd([char]13, [char]7)
        if ($nextText -eq "" -and $next.Range.ParagraphStyle.NameLocal -eq "No Spacing") {
            $next.Range.Delete()
        }
        break
    }
}
